$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E86").Value = -3.101195725775879
$ws.Range("F86").Value = 3.720222109573349
$ws.Range("I86").Value = -2.934478420543018

$ws.Range("E87").Value = 3.932030374331585
$ws.Range("F87").Value = 4.681962787741607
$ws.Range("I87").Value = 1.785321559763111

$ws.Range("E88").Value = -2.069715247260751
$ws.Range("F88").Value = 0.5920895797350152
$ws.Range("I88").Value = -1.631227848023713

$ws.Range("E89").Value = -4.698336107940325
$ws.Range("F89").Value = -1.484304176661342
$ws.Range("I89").Value = -3.670248356208035

$ws.Range("E90").Value = -2.286601610844267
$ws.Range("F90").Value = -1.280655647928439
$ws.Range("I90").Value = -2.150384566497852

$ws.Range("E91").Value = -2.339571895058551
$ws.Range("F91").Value = -2.848556215275973
$ws.Range("I91").Value = -1.732726665771843

$ws.Range("E92").Value = -2.13930466170003
$ws.Range("F92").Value = -2.865953568885793
$ws.Range("I92").Value = -2.070406934421564

$ws.Range("E93").Value = -0.7357084438719351
$ws.Range("F93").Value = -1.875296652868696
$ws.Range("I93").Value = -0.6287371506267943
